# Applies:
#  1. The table on slide 16 switches its table style (tableStyleId) from
#     {436F9605-DD0E-4399-91F9-4696435E9B8D} to {BD3C93A1-9E8B-4F9E-8EC6-B2F39B0B4E06}.
#  2. The presentation's theme colour scheme (ppt/theme/theme1.xml, the theme
#     used by the slide master / all slides) is repainted from the "Integral"
#     palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$oldStyleId = "{436F9605-DD0E-4399-91F9-4696435E9B8D}"
$newStyleId = "{BD3C93A1-9E8B-4F9E-8EC6-B2F39B0B4E06}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.StyleId -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Theme colours -------------------------------------------------------
# Order exposed by ThemeColorScheme.Item(n): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
